$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Move the existing "Week 2" task column (C17:C24) into a new
#    "Week 3" column (D17:D24), then fill in the new C (Dataset
#    Discovery) and E (Week 3 attribute-analysis) columns.
# -----------------------------------------------------------------
$oldC = @{}
for ($r = 17; $r -le 24; $r++) {
    $oldC[$r] = $ws.Cells.Item($r, 3).Value()
}
for ($r = 17; $r -le 24; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldC[$r]
}

$ws.Range("C17").Value = "Dataset Discovery"
$ws.Range("C18").Value = "Dataset Discovery"
$ws.Range("C19").Value = "Dataset Discovery"
$ws.Range("C20").Value = "Dataset Discovery"
$ws.Range("C21").Value = "Dataset Discovery"
$ws.Range("C22").Value = "Dataset Discovery"
$ws.Range("C23").Value = "Dataset Discovery"
$ws.Range("C24").Value = "Dataset Discovery"

$ws.Range("E17").Value = "Presentation preparation & attribute analysis(Partner,Dependents,tenure)"
$ws.Range("E18").Value = "Attribute introduction,analysis and presentation preparation"
$ws.Range("E19").Value = "Attribute analysis(gender & seniorCitizen)"
$ws.Range("E20").Value = "Attribute analysis(Contract & PaperlessBilling)"
$ws.Range("E21").Value = "Attribute analysis(DeviceProtection & TechSupport)"
$ws.Range("E22").Value = "Attribute analysis(OnlineSecurity,OnlineBackup)"
$ws.Range("E23").Value = "Presentation Preparation & Dataset understanding(PhoneService,MultipleLines,InternetService)"
$ws.Range("E24").Value = "Dataset analysis and visualization(StreamingTV,StreamingMovies)"

# -----------------------------------------------------------------
# 2. Update the weekly "Notes" row (row 28): reword the Week-2 note,
#    and add Week-3 / Week-4 notes in the new D28 / E28 cells.
# -----------------------------------------------------------------
$ws.Range("C28").Value = "On week 2, The project was decided on and the team was introduced to the dataset and its attributes as well as the introduction to kaggle."
$ws.Range("D28").Value = "On Week 3, Each member was introduced to the dataset,GitHub, Kaggle and Jupyter notebook."
$ws.Range("E28").Value = "Week 4, Dataset was divided among each member and was introduced to the basics of Jupyter notebook, Pandas and matplotlib. Presentation was created by the team with everybody present, so that we could brainstorm and collaborate to get suggestions."

# Formatting for row 28: wrap text + top vertical alignment across the
# whole note row; A28 keeps its green fill, E28 keeps the extra
# "horizontal-left" alignment that the rest of the row doesn't have.
$ws.Range("A28:E28").WrapText = $true
$ws.Range("A28:E28").VerticalAlignment = -4160
$ws.Range("E28").HorizontalAlignment = -4131
$ws.Rows.Item(28).RowHeight = 57.6

# -----------------------------------------------------------------
# 3. Column sizing for the now-used B:E columns.
# -----------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.608072916666668
$ws.Columns.Item(2).ColumnWidth = 37.053385416666664
$ws.Columns.Item(3).ColumnWidth = 37.276041666666664
$ws.Columns.Item(4).ColumnWidth = 37.053385416666664
$ws.Columns.Item(5).ColumnWidth = 78.83072916666667

# -----------------------------------------------------------------
# 4. Misc leftover cell (H16) that extends the sheet's used range out
#    to column H, matching the refreshed workbook dimension.
# -----------------------------------------------------------------
$ws.Range("H16").Font.Bold = $false

# -----------------------------------------------------------------
# 5. View state: scroll down a bit and select C29:C30 like the
#    author's last saved session.
# -----------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("C29:C30").Select()
